# Fix subscript issue in spreadsheet:
#  - Replace subscript/odd header strings (k , E1..E5/M1..M5 with Unicode
#    subscripts) with plain ASCII header strings across every sheet, so the
#    shared-string table collapses down to the 11 already-existing plain
#    variants (k, E1, E2, E3, E4, E5, M1, M2, M3, M4, M5).
#  - Reset most sheets' saved selection back to A1 (sheet "72" keeps A2).
#  - Make sheet "92" the active sheet/tab, with its selection at M18.
#
# NOTE: hashtable (@{...}) iteration via foreach(...Keys) proved unreliable
# in this COM-interop host (intermittent "subscript out of range" / null
# method calls), so plain parallel arrays + indexed for-loops are used
# instead everywhere below.

$wb = $excel.ActiveWorkbook

$headers = @("k", "E1", "E2", "E3", "E4", "E5", "M1", "M2", "M3", "M4", "M5")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    for ($col = 1; $col -le $headers.Length; $col++) {
        $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
    }
}

# Per-sheet saved selection (by tab name, matching the diff). Sheet "72" is
# intentionally omitted - it keeps its original A2 selection.
$selNames = @("10", "20", "30", "40", "54", "64", "78", "83", "88")
$selCells = @("A1", "A1", "A1", "A1", "A1", "A1", "A1", "A1", "A1")

for ($i = 0; $i -lt $selNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($selNames[$i])
    $ws.Range($selCells[$i]).Select()
}

# Activate sheet "92" last so it becomes the active tab/sheet, and its
# selection (M18) sticks as the active-window selection.
$target = $wb.Worksheets.Item("92")
$target.Activate()
$target.Range("M18").Select()
